$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.531.86'
$ws.Range("E2").Value = '  +1.43%  '

$ws.Range("D3").Value = '3.451.29'
$ws.Range("E3").Value = '  +2.43%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.85%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.94'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.97%  '

$ws.Range("D7").Value = '3.451.07'
$ws.Range("E7").Value = '  +2.51%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.64'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.79%  '

$ws.Range("E11").Value = '  +3.76%  '

$ws.Range("E12").Value = '  +2.33%  '

$ws.Range("D13").Value = '4.043.87'
$ws.Range("E13").Value = '  +2.51%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.08'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +8.72%  '

$ws.Range("E15").Value = '  -1.08%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000173'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.14%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.450.58'
$ws.Range("E17").Value = '  +2.37%  '

$ws.Range("D18").Value = '61.704.87'
$ws.Range("E18").Value = '  +1.47%  '

$ws.Range("E19").Value = '  +7.79%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.26'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.42%  '

$ws.Range("E21").Value = '  +3.40%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '391.03'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.09%  '

$ws.Range("E23").Value = '  +3.53%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.30'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.51%  '

$ws.Range("E25").Value = '  -0.06%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.75'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.07%  '

$ws.Range("E27").Value = '  +0.88%  '

$ws.Range("D28").Value = '3.591.25'
$ws.Range("E28").Value = '  +2.31%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.179'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.97%  '

$ws.Range("E30").Value = '  +4.81%  '

$ws.Range("E31").Value = '  +0.12%  '

$ws.Range("E32").Value = '  -8.98%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.13'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.21%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '24.04'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.34%  '

$ws.Range("D37").Value = '3.481.86'
$ws.Range("E37").Value = '  +2.62%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.01'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.77%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.12'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.60%  '

$ws.Range("E40").Value = '  +1.78%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '167.18'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.81%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '28.34'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +13.16%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0781'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.67%  '

$ws.Range("E44").Value = '  +4.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.73'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.73%  '

$ws.Range("E46").Value = '  +0.04%  '

$ws.Range("E47").Value = '  +4.48%  '

$ws.Range("E48").Value = '  +1.54%  '

$ws.Range("D49").Value = '2.592.07'
$ws.Range("E49").Value = '  +3.04%  '

$ws.Range("E50").Value = '  -0.63%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.91'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.52%  '
